$d = $word.ActiveDocument

# The document currently ends with:
#   ... "Answers to questions *before* writing, or running, code:"
#   <empty paragraph>
# Two new numbered-list items need to be inserted between those two
# paragraphs.

$intro = $d.Paragraphs.Item(2)

# Insert the first new list item right after the intro paragraph.
$intro.Range.InsertParagraphAfter()
$item1 = $d.Paragraphs.Item(3)
$item1.Range.Text = "The theoretical Big-O execution time *should* be linear but, as we learned in class, our code won’t because of the implementation."

# Insert the second new list item right after the first one.
$item1.Range.InsertParagraphAfter()
$item2 = $d.Paragraphs.Item(4)
$item2.Range.Text = "The internet says the same thing – that it should be linear."

# Apply the "List Paragraph" style and default numbering (1., 2., ...) to
# both new paragraphs together so they belong to the same numbered list.
$listRange = $d.Range($item1.Range.Start, $item2.Range.End)
$listRange.Style = "List Paragraph"
$listRange.ListFormat.ApplyNumberDefault()
